# Update countries & provincias Spain
#
# The source workbook lists one country per row (columns B-H hold that
# day's case counts). The new day's export re-sorted a handful of
# countries, which - because the sheet stores a shared-string index per
# row - shows up as the country name *and* its numeric row changing for
# the rows whose sort position moved. The "Datos actualizados" timestamp
# (cell A1) also moved from 21:03 to 22:03.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> index
$col = @{ "A"=1; "B"=2; "C"=3; "D"=4; "E"=5; "F"=6; "G"=7; "H"=8 }

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 22:03"

# Per-row updates: country name (where it changed) + new totals.
$updates = @{
  4   = @{ B=1229919; C=17084; D=198935; E=959314; F=16174; G=1749; H=71670 }
  9   = @{ B=166696; C=544; E=24603 }
  16  = @{ D=15413; E=34332 }

  49  = @{ A="Sudafrica"; B=7572; C=352;  D=2746; E=4678; F=36; G=10; H=148 }
  50  = @{ A="Panama";    B=7387; C=190;  D=726;  E=6458; F=93; G=3;  H=203 }

  62  = @{ B=3720; C=187; E=1950 }

  77  = @{ A="Guinea";   B=1811; C=101; D=498;  E=1303; G=1 }
  78  = @{ A="Islandia"; B=1799; C=0;   D=1733; E=56;   F=0; H=10 }
  79  = @{ A="Estonia";  B=1711; C=8;   D=261;  E=1395; F=6; H=55 }

  101 = @{ A="Costa Rica";             C=13; D=413; E=336; F=5;  H=6 }
  102 = @{ A="Niger";                  B=755; C=0;  D=534; E=184; F=0;  G=0; H=37 }
  103 = @{ A="Principado de Andorra";  B=751; C=1;  D=514; E=191; F=16; G=1; H=46 }

  115 = @{ D=20; E=551; G=1; H=2 }

  125 = @{ A="Estado de Palestina"; B=371; C=9; D=127; E=242; F=0; H=2 }
  126 = @{ A="Gabon";               B=367; C=0; D=93;  E=268; F=1; H=6 }

  166 = @{ B=81; C=1; E=62 }

  199 = @{ A="Burundi";                D=7; H=1 }
  200 = @{ A="San Cristobal y Nieves"; D=8; H=0 }

  205 = @{ D=8; E=3 }
}

foreach ($r in $updates.Keys) {
  $rowVals = $updates[$r]
  foreach ($c in $rowVals.Keys) {
    $ws.Cells.Item([int]$r, $col[$c]).Value = $rowVals[$c]
  }
}
